# Add a blank (spacer) paragraph right after the last table and before the
# "Plus de projets..." paragraph, matching the formatting (tabs, spacing,
# italic run-properties on the paragraph mark) used elsewhere in the doc.

$d = $word.ActiveDocument

# Locate the last table in the document and the position right after it,
# i.e. the start of the paragraph that currently follows the table.
$lastTable = $d.Tables.Item($d.Tables.Count)
$insertPos = $lastTable.Range.End

$target = $d.Range($insertPos, $insertPos)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:pPr>' + `
               '<w:tabs>' + `
                   '<w:tab w:val="left" w:pos="2160"/>' + `
                   '<w:tab w:val="left" w:pos="2880"/>' + `
               '</w:tabs>' + `
               '<w:spacing w:before="160" w:after="0"/>' + `
               '<w:rPr>' + `
                   '<w:i/>' + `
                   '<w:iCs/>' + `
                   '<w:lang w:val="fr-FR"/>' + `
               '</w:rPr>' + `
           '</w:pPr>' + `
       '</w:p>'

$target.InsertXML($xml)
